$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H76").Value = 275002140
$ws.Range("I76").Value = 275002140
$ws.Range("K76").Value = 275002140
$ws.Range("M76").Value = -275001825
$ws.Range("H79").Value = 275002140
$ws.Range("I79").Value = 275002140
$ws.Range("K79").Value = 275002140
$ws.Range("M79").Value = -275001048
$ws.Range("H106").Value = 2934.182
$ws.Range("I106").Value = 3063.3333
$ws.Range("J106").Value = 2353
$ws.Range("K106").Value = 3063.3333
$ws.Range("L106").Value = 2353
$ws.Range("M106").Value = -2432.3333
$ws.Range("N106").Value = -3615
$ws.Range("H135").Value = 1763.6316
$ws.Range("I135").Value = 853.4706
$ws.Range("J135").Value = 9500
$ws.Range("K135").Value = 7681.2354
$ws.Range("L135").Value = 85500
$ws.Range("M135").Value = -5146.2354
$ws.Range("N135").Value = -90570
$ws.Range("H137").Value = 1200
$ws.Range("I137").Value = 1200
$ws.Range("K137").Value = 3600
$ws.Range("M137").Value = -1050

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 1631658
$ws.Range("I32").Value = 4290.2646
$ws.Range("K32").Value = 4290.2646
$ws.Range("M32").Value = -4003.2646
$ws.Range("H61").Value = 1258.7778
$ws.Range("I61").Value = 1168.7632
$ws.Range("J61").Value = 1472.5625
$ws.Range("K61").Value = 1168.7632
$ws.Range("L61").Value = 1472.5625
$ws.Range("M61").Value = -956.7632000000001
$ws.Range("N61").Value = -1896.5625
$ws.Range("H97").Value = 541.52
$ws.Range("I97").Value = 529.94446
$ws.Range("J97").Value = 571.2857
$ws.Range("K97").Value = 529.94446
$ws.Range("L97").Value = 571.2857
$ws.Range("M97").Value = -33.94446000000005
$ws.Range("N97").Value = -1563.2857
$ws.Range("H136").Value = 1258.7778
$ws.Range("I136").Value = 1168.7632
$ws.Range("J136").Value = 1472.5625
$ws.Range("K136").Value = 3506.2896
$ws.Range("L136").Value = 4417.6875
$ws.Range("M136").Value = -956.2896000000001
$ws.Range("N136").Value = -9517.6875
$ws.Range("H139").Value = 44707.5
$ws.Range("J139").Value = 44707.5
$ws.Range("L139").Value = 44707.5
$ws.Range("N139").Value = -54987.5

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H94").Value = 733.7727
$ws.Range("I94").Value = 590.6923
$ws.Range("J94").Value = 940.44446
$ws.Range("K94").Value = 590.6923
$ws.Range("L94").Value = 940.44446
$ws.Range("M94").Value = -139.6923
$ws.Range("N94").Value = -1842.44446
$ws.Range("H137").Value = 51819
$ws.Range("J137").Value = 51819
$ws.Range("L137").Value = 51819
$ws.Range("N137").Value = -62019

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H62").Value = 10492.308
$ws.Range("J62").Value = 8800
$ws.Range("L62").Value = 8800
$ws.Range("N62").Value = -10048
$ws.Range("H65").Value = 10492.308
$ws.Range("J65").Value = 8800
$ws.Range("L65").Value = 44000
$ws.Range("N65").Value = -50240
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H134").Value = 1605.025
$ws.Range("I134").Value = 1338.9375
$ws.Range("K134").Value = 4016.8125
$ws.Range("M134").Value = -1481.8125
$ws.Range("H140").Value = 49403.332
$ws.Range("J140").Value = 49403.332
$ws.Range("L140").Value = 49403.332
$ws.Range("N140").Value = -59763.332

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H63").Value = 87799.914
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 87799.914
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 263399.742
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -264897.742
$ws.Range("H66").Value = 87799.914
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 87799.914
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 790199.226
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -797687.226
$ws.Range("H68").Value = 668
$ws.Range("I68").Value = 661.6
$ws.Range("J68").Value = 700
$ws.Range("K68").Value = 1984.8
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -1173.8
$ws.Range("N68").Value = -3722
$ws.Range("H71").Value = 668
$ws.Range("I71").Value = 661.6
$ws.Range("J71").Value = 700
$ws.Range("K71").Value = 5954.400000000001
$ws.Range("L71").Value = 6300
$ws.Range("M71").Value = -1898.400000000001
$ws.Range("N71").Value = -14412
$ws.Range("H132").Value = 860
$ws.Range("I132").Value = 755
$ws.Range("K132").Value = 6795
$ws.Range("M132").Value = -4265

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 4381.778
$ws.Range("I70").Value = 4136.32
$ws.Range("J70").Value = 7450
$ws.Range("K70").Value = 4136.32
$ws.Range("L70").Value = 7450
$ws.Range("M70").Value = -3866.32
$ws.Range("N70").Value = -7990
$ws.Range("H73").Value = 4381.778
$ws.Range("I73").Value = 4136.32
$ws.Range("J73").Value = 7450
$ws.Range("K73").Value = 4136.32
$ws.Range("L73").Value = 7450
$ws.Range("M73").Value = -3200.32
$ws.Range("N73").Value = -9322
$ws.Range("H80").Value = 3499.9412
$ws.Range("I80").Value = 3513.2666
$ws.Range("K80").Value = 3513.2666
$ws.Range("M80").Value = -2515.2666
$ws.Range("H83").Value = 3499.9412
$ws.Range("I83").Value = 3513.2666
$ws.Range("K83").Value = 17566.333
$ws.Range("M83").Value = -12574.333
$ws.Range("H113").Value = 1502.875
$ws.Range("I113").Value = 1335
$ws.Range("J113").Value = 2006.5
$ws.Range("K113").Value = 1335
$ws.Range("L113").Value = 2006.5
$ws.Range("M113").Value = 835
$ws.Range("N113").Value = -6346.5
$ws.Range("H123").Value = 25739.5
$ws.Range("J123").Value = 25739.5
$ws.Range("L123").Value = 25739.5
$ws.Range("N123").Value = -30639.5
$ws.Range("H135").Value = 29800
$ws.Range("J135").Value = 29800
$ws.Range("L135").Value = 29800
$ws.Range("N135").Value = -39940
$ws.Range("H138").Value = 35857.25
$ws.Range("J138").Value = 35857.25
$ws.Range("L138").Value = 35857.25
$ws.Range("N138").Value = -46137.25

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 3134
$ws.Range("I61").Value = 2600.8572
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2600.8572
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2398.8572
$ws.Range("N61").Value = -5404
$ws.Range("H100").Value = 5908142.5
$ws.Range("I100").Value = 6602971
$ws.Range("K100").Value = 6602971
$ws.Range("M100").Value = -6602430
$ws.Range("H113").Value = 3134
$ws.Range("I113").Value = 2600.8572
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2600.8572
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -430.8571999999999
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 6407.684
$ws.Range("I122").Value = 18070.4
$ws.Range("J122").Value = 2242.4285
$ws.Range("K122").Value = 54211.2
$ws.Range("L122").Value = 6727.2855
$ws.Range("M122").Value = -51761.2
$ws.Range("N122").Value = -11627.2855
$ws.Range("H134").Value = 47425.8
$ws.Range("J134").Value = 47425.8
$ws.Range("L134").Value = 47425.8
$ws.Range("N134").Value = -57565.8
$ws.Range("H136").Value = 2531.2415
$ws.Range("I136").Value = 1183.5217
$ws.Range("J136").Value = 7697.5
$ws.Range("K136").Value = 3550.5651
$ws.Range("L136").Value = 23092.5
$ws.Range("M136").Value = -1000.5651
$ws.Range("N136").Value = -28192.5
$ws.Range("H138").Value = 35643
$ws.Range("J138").Value = 35643
$ws.Range("L138").Value = 35643
$ws.Range("N138").Value = -45923

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H7").Value = 23669.834
$ws.Range("I7").Value = 502
$ws.Range("J7").Value = 35253.75
$ws.Range("K7").Value = 502
$ws.Range("L7").Value = 35253.75
$ws.Range("M7").Value = -389
$ws.Range("N7").Value = -35479.75
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10586
$ws.Range("H55").Value = 12250
$ws.Range("I55").Value = 4500
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -4223
$ws.Range("N55").Value = -20554
$ws.Range("H100").Value = 780.13336
$ws.Range("I100").Value = 256.14285
$ws.Range("J100").Value = 1238.625
$ws.Range("K100").Value = 512.2857
$ws.Range("L100").Value = 2477.25
$ws.Range("M100").Value = 28.71429999999998
$ws.Range("N100").Value = -3559.25
$ws.Range("H122").Value = 1427.4286
$ws.Range("I122").Value = 1700.6666
$ws.Range("K122").Value = 5101.9998
$ws.Range("M122").Value = -2651.9998
$ws.Range("H136").Value = 1716.1111
$ws.Range("I136").Value = 1647.5
$ws.Range("K136").Value = 4942.5
$ws.Range("M136").Value = -2392.5
